$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws1.Activate()
$excel.Goto($ws1.Range("A97"), $true)
$ws1.Range("J118").Select()
Write-Host "done"
